# "leaves code up to date till 6/28/2019"
#
# On the "Config" sheet, the TestCases value (cell B2) is updated from a
# single, stale test-case id ("16") to the up-to-date, comma-separated
# list of ids that should now run ("26,28,29,30,31,36").
#
# Excel auto-sizes the (now much longer) column B to fit the new value,
# so we widen it to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Config")

$ws.Range("B2").Value = "26,28,29,30,31,36"

# Best-fit column B to the new, longer text (~13.3 chars wide).
$ws.Columns("B:B").ColumnWidth = 12.5
